$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# B8 / B9: planned calendar week (written first so the new shared string
# "KW 33" lands at index 18, before "a" below, matching the target file)
$ws.Range("B8").Value = "KW 33"
$ws.Range("B9").Value = "KW 33"

# Row 5: new value "a" in C5 with a centered Marlett 12pt font, and a
# slightly taller row to fit it.
$ws.Range("C5").Value = "a"
$ws.Range("C5").Font.Name = "Marlett"
$ws.Range("C5").Font.Size = 12
$ws.Range("C5").HorizontalAlignment = -4108
$ws.Rows.Item(5).RowHeight = 17.25

# Move the active selection to E9
$ws.Range("E9").Select()
